$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 (the "be the best to make some plank..." review row),
# which shifts all the rows below it up by one.
$ws.Rows.Item(3).Delete()

# The hyperlinks that were anchored to the old row positions need to be
# re-pointed to the new (shifted-up) cell locations. Clear existing
# hyperlinks and recreate them at the correct cells.
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:stevewonder3001@gmail.com", "", "", "stevewonder3001@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:budoyoni@gmail.com", "", "", "budoyoni@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:redvelvetmichael@gmail.com", "", "", "redvelvetmichael@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:veredsnir12@gmail.com", "", "", "veredsnir12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:vikicrestina@gmail.com", "", "", "vikicrestina@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:cristianjohn1222@gmail.com", "", "", "cristianjohn1222@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:bittonnir12@gmail.com", "", "", "bittonnir12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:nevilgreen12@gmail.com", "", "", "nevilgreen12@gmail.com")

# Restore the selection to where the user last clicked after the edit.
$ws.Range("B3").Select()
